# Add B3 R3 leachate mass data.
#
# Fills in the "water_plus_bottle_mass_collected(g)" (F) and
# "DNA_filter_date" (H) columns for the B3 R3 samples (rows 38-55 on the
# "Rainfall 2" sheet). The existing "water_mass(g)" column (G) is a shared
# formula (=F-E) that recalculates automatically once F is populated.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rainfall 2")
$ws.Activate()

$filterDate = 20210611

$waterPlusBottleMass = @{
    38 = 772
    39 = 814
    40 = 905
    41 = 843
    42 = 827
    43 = 810
    44 = 847
    45 = 790
    46 = 857
    47 = 833
    48 = 825
    49 = 869
    50 = 868
    51 = 845
    52 = 855
    53 = 854
    54 = 870
    55 = 861
}

foreach ($row in 38..55) {
    $ws.Cells.Item($row, 6).Value = $waterPlusBottleMass[$row]
    $ws.Cells.Item($row, 8).Value = $filterDate
}

# Leave the view positioned where the author ended up after the edit.
$ws.Range("F56").Select()
